$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.528161333333333
$ws.Range("H2").Value = 13.584484
$ws.Range("I2").Value = 0.3225352762763812
$ws.Range("J2").Value = 0.3225352762763812
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.35607266666667
$ws.Range("N2").Value = 37.068218
$ws.Range("Q2").Value = 55.9502904810569
$ws.Range("R2").Value = 503.552614329512
$ws.Range("S2").Value = 0.3225352762763812
$ws.Range("T2").Value = 0.3225352762763812

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.270951666666666
$ws.Range("H3").Value = 18.812855
$ws.Range("I3").Value = 0.4466720550425397
$ws.Range("J3").Value = 0.4466720550425397
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.35607266666667
$ws.Range("N3").Value = 37.068218
$ws.Range("Q3").Value = 77.48433448248778
$ws.Range("R3").Value = 697.35901034239
$ws.Range("S3").Value = 0.4466720550425397
$ws.Range("T3").Value = 0.4466720550425397

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.240161666666667
$ws.Range("H4").Value = 9.720485
$ws.Range("I4").Value = 0.2307926686810791
$ws.Range("J4").Value = 0.2307926686810791
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.35607266666667
$ws.Range("N4").Value = 37.068218
$ws.Range("Q4").Value = 40.03567300508112
$ws.Range("R4").Value = 360.32105704573
$ws.Range("S4").Value = 0.2307926686810791
$ws.Range("T4").Value = 0.2307926686810791
